$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cell, $value)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "271.86"
Set-TextValue $ws.Range("G2") "17"
Set-TextValue $ws.Range("D3") "23.13"
Set-TextValue $ws.Range("G3") "17"
Set-TextValue $ws.Range("D4") "6.362"
Set-TextValue $ws.Range("G4") "17"
Set-TextValue $ws.Range("D5") "0.06302"
Set-TextValue $ws.Range("G5") "17"
Set-TextValue $ws.Range("D6") "3.658"
Set-TextValue $ws.Range("G6") "17"
Set-TextValue $ws.Range("D7") "6.766"
Set-TextValue $ws.Range("G7") "17"
Set-TextValue $ws.Range("D8") "1.393"
Set-TextValue $ws.Range("G8") "17"
Set-TextValue $ws.Range("D9") "0.8339"
Set-TextValue $ws.Range("G9") "17"
Set-TextValue $ws.Range("D10") "0.1629"
Set-TextValue $ws.Range("G10") "17"
Set-TextValue $ws.Range("D11") "0.08474"
Set-TextValue $ws.Range("G11") "17"
Set-TextValue $ws.Range("D12") "0.03469"
Set-TextValue $ws.Range("G12") "17"
Set-TextValue $ws.Range("D13") "0.03125"
Set-TextValue $ws.Range("G13") "17"
Set-TextValue $ws.Range("D14") "0.09311"
Set-TextValue $ws.Range("G14") "17"
Set-TextValue $ws.Range("D15") "3.940"
Set-TextValue $ws.Range("G15") "17"
Set-TextValue $ws.Range("D16") "0.001700"
Set-TextValue $ws.Range("G16") "17"
Set-TextValue $ws.Range("D17") "0.04863"
Set-TextValue $ws.Range("G17") "17"
Set-TextValue $ws.Range("D18") "0.006210"
Set-TextValue $ws.Range("G18") "17"
Set-TextValue $ws.Range("D19") "0.005480"
Set-TextValue $ws.Range("G19") "17"
Set-TextValue $ws.Range("D20") "0.001089"
Set-TextValue $ws.Range("G20") "17"
Set-TextValue $ws.Range("D21") "0.0001499"
Set-TextValue $ws.Range("G21") "17"
Set-TextValue $ws.Range("D22") "3.731"
Set-TextValue $ws.Range("G22") "17"
Set-TextValue $ws.Range("D23") "2.355"
Set-TextValue $ws.Range("G23") "17"
Set-TextValue $ws.Range("D24") "0.01387"
Set-TextValue $ws.Range("G24") "17"
Set-TextValue $ws.Range("G25") "17"
Set-TextValue $ws.Range("G26") "17"
Set-TextValue $ws.Range("D27") "0.0003740"
Set-TextValue $ws.Range("G27") "17"
Set-TextValue $ws.Range("G28") "17"
Set-TextValue $ws.Range("G29") "17"
Set-TextValue $ws.Range("G30") "17"
Set-TextValue $ws.Range("G31") "17"
Set-TextValue $ws.Range("G32") "17"
Set-TextValue $ws.Range("G33") "17"
Set-TextValue $ws.Range("G34") "17"
Set-TextValue $ws.Range("G35") "17"
Set-TextValue $ws.Range("G36") "17"
Set-TextValue $ws.Range("G37") "17"
Set-TextValue $ws.Range("G38") "17"
Set-TextValue $ws.Range("G39") "17"
Set-TextValue $ws.Range("D40") "0.04692"
Set-TextValue $ws.Range("G40") "17"
Set-TextValue $ws.Range("D41") "0.006897"
Set-TextValue $ws.Range("G41") "17"
Set-TextValue $ws.Range("G42") "17"
Set-TextValue $ws.Range("D43") "0.003597"
Set-TextValue $ws.Range("G43") "17"
Set-TextValue $ws.Range("D44") "0.01251"
Set-TextValue $ws.Range("G44") "17"
Set-TextValue $ws.Range("D45") "0.00006244"
Set-TextValue $ws.Range("G45") "17"
Set-TextValue $ws.Range("D46") "0.00000000749"
Set-TextValue $ws.Range("G46") "17"
Set-TextValue $ws.Range("D47") "0.7995"
$ws.Range("E47").Value = "46CoinbaseStockTokenCOIN"
Set-TextValue $ws.Range("G47") "17"
Set-TextValue $ws.Range("D48") "0.1114"
Set-TextValue $ws.Range("G48") "17"
Set-TextValue $ws.Range("D49") "0.00002099"
$ws.Range("E49").Value = "48CryptobidCoinCBCWorstin24h"
Set-TextValue $ws.Range("G49") "17"
Set-TextValue $ws.Range("D50") "0.01239"
Set-TextValue $ws.Range("G50") "17"
Set-TextValue $ws.Range("G51") "17"
